# Add the author's own "tache" (task) rows to the Tableau_Taches table.
# Column layout: A=Etape de modelisation, B=Projet, C=Tache,
#                D=Description, E=Date debut, F=Date Fin

$xlVAlignBottom = -4107
$xlHAlignCenter = -4108

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").VerticalAlignment = $xlVAlignBottom
$ws.Range("A2").Value = 'Analyse fonctionelle'
$ws.Range("B2").Value = 'Vtts-windows-app'
$ws.Range("C2").VerticalAlignment = $xlVAlignBottom
$ws.Range("C2").Value = 'Diagramme ca utilisation Securitie gwin'
$ws.Range("D2").Value = 'Diagramme cas utilisation Securitie pour  Gwin'
$ws.Range("E2").Value = 42863
$ws.Range("F2").Value = 42863

# --- Row 3 ---
$ws.Range("A3").VerticalAlignment = $xlVAlignBottom
$ws.Range("A3").Value = 'Analyse fonctionelle'
$ws.Range("B3").Value = 'Sport-club-management'
$ws.Range("C3").VerticalAlignment = $xlVAlignBottom
$ws.Range("C3").Value = 'Diagramme de classe  '
$ws.Range("D3").Value = 'Diagramme de classz de sport clob management'
$ws.Range("E3").Value = 42853
$ws.Range("F3").Value = 42853

# --- Row 4 ---
$ws.Range("A4").VerticalAlignment = $xlVAlignBottom
$ws.Range("A4").Value = 'Réalisation '
$ws.Range("B4").Value = 'Sport-club-management'
$ws.Range("C4").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C4").Value = 'Add All class to project c#'
$ws.Range("D4").Value = 'Ajouter tous les class au projet avec c#'
$ws.Range("E4").Value = 42857
$ws.Range("F4").Value = 42858

# --- Row 5 ---
$ws.Range("A5").VerticalAlignment = $xlVAlignBottom
$ws.Range("A5").Value = 'Réalisation '
$ws.Range("B5").Value = 'Sport-club-management'
$ws.Range("C5").VerticalAlignment = $xlVAlignBottom
$ws.Range("C5").Value = 'Securitie'
$ws.Range("D5").Value = 'ajouter permission au utilisateur concerné'
$ws.Range("E5").Value = 42857
$ws.Range("F5").Value = 42857

# --- Row 6 ---
$ws.Range("A6").VerticalAlignment = $xlVAlignBottom
$ws.Range("A6").Value = 'Réalisation '
$ws.Range("B6").Value = 'Sport-club-management'
$ws.Range("C6").Value = 'compléter la tradiction des propriétés'
$ws.Range("D6").Value = 'compléter la tradiction des propriétés dans une fichier XL'
$ws.Range("E6").Value = 42850
$ws.Range("F6").Value = 42850

# --- Row 7 ---
$ws.Range("A7").VerticalAlignment = $xlVAlignBottom
$ws.Range("A7").Value = 'Réalisation '
$ws.Range("B7").Value = 'Sport-club-management'
$ws.Range("C7").Value = 'Modifier les type string'
$ws.Range("D7").Value = 'modifier tous les variable de type string en type LocalizedString'
$ws.Range("E7").Value = 42860
$ws.Range("F7").Value = 42860

# --- Row 8 ---
$ws.Range("A8").Value = 'Formation'
$ws.Range("B8").Value = 'Sport-club-management'
$ws.Range("C8").Value = 'formation sur github'
$ws.Range("D8").Value = 'formation sur github'
$ws.Range("E8").Value = 42857
$ws.Range("F8").Value = 42857

# --- Row 9 ---
$ws.Range("A9").Value = 'Formation'
$ws.Range("B9").Value = 'Sport-club-management'
$ws.Range("C9").Value = 'informer stagiaire sur github'
$ws.Range("D9").Value = 'informer stagiaire sur github'
$ws.Range("E9").Value = 42859
$ws.Range("F9").Value = 42859

# --- Row 10 ---
$ws.Range("A10").Value = 'REALISATION'
$ws.Range("B10").Value = 'Vtts-windows-app'
$ws.Range("C10").Value = 'Ajouter des autorisation'
$ws.Range("D10").Value = 'ajouter autorisation pour  utilisateur admin '
$ws.Range("E10").Value = 42863
$ws.Range("F10").Value = 42863

# --- Row 11 ---
$ws.Range("A11").Value = 'Realisation'
$ws.Range("B11").Value = 'Vtts-windows-app'
$ws.Range("C11").Value = 'ajouter les type des specialities'
$ws.Range("D11").Value = 'ajouter les types des specialitie par defeaut'
$ws.Range("E11").Value = 42863
$ws.Range("F11").Value = 42863

# --- Row 12 ---
$ws.Range("A12").Value = 'REALISATION'
$ws.Range("B12").Value = 'Vtts-windows-app'
$ws.Range("C12").Value = 'dimonstration vedio github'
$ws.Range("D12").Value = 'dimonstration vedio Pour expliquer comment la facon de update projet'
$ws.Range("E12").Value = 42863
$ws.Range("F12").Value = 42863

# The wrapped text needs two or three display lines in several rows;
# reflect that in the row heights (mirrors what Excel's own autofit
# would compute once the text was entered).
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 45

# Restore the active selection to match the edited range
$ws.Range("A2:F12").Select()
